$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; this shifts the existing data rows (2..20)
# down to (3..21) and carries their values/styles along, but this host's
# Hyperlinks collection does not re-target itself to the shifted rows, so
# we rebuild the hyperlinks list from scratch further down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the latest circular entry.
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 290.95
$ws.Range("E2").Value = "27-11-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-27-11-2025.pdf"

# The existing hyperlinks collection became misaligned after the row insert
# (their anchors did not shift down with the cells). Clear it out and
# rebuild every hyperlink, in row order, against the final layout.
$ws.Range("F2").Hyperlinks.Delete()

$links = @(
  @{ Row = 2;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-27-11-2025.pdf" },
  @{ Row = 3;  Url = "https://nalcoindia.com/wp-content/uploads/2025/11/Ingot-22-11-2025.pdf" },
  @{ Row = 4;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf" },
  @{ Row = 5;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf" },
  @{ Row = 6;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf" },
  @{ Row = 7;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf" },
  @{ Row = 8;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf" },
  @{ Row = 9;  Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf" },
  @{ Row = 10; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf" },
  @{ Row = 11; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf" },
  @{ Row = 12; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf" },
  @{ Row = 13; Url = "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf" },
  @{ Row = 14; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf" },
  @{ Row = 15; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf" },
  @{ Row = 16; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf" },
  @{ Row = 17; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf" },
  @{ Row = 18; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf" },
  @{ Row = 19; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf" },
  @{ Row = 20; Url = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf" },
  @{ Row = 21; Url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf" }
)

foreach ($link in $links) {
  $cell = $ws.Cells.Item($link.Row, 6)
  $ws.Hyperlinks.Add($cell, $link.Url, "", "", $link.Url)
}
